# GestioneSalariOspedale.xlsx - "Modifiche varie nomi e templates"
#
# The G4 cell used to hold a broken formula ( D4 +( E4 * F4 ) ) that always
# evaluated to #VALUE! because F4 contains text (the "HourlyPay" template
# placeholder), not a number. The fix replaces that formula with the literal
# template placeholder string so the export engine can substitute it later,
# which also makes the SUM() in G5 resolve cleanly instead of erroring out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the broken "=D4+(E4*F4)" formula in G4 with the literal
# placeholder text. Clear the formula first, then set the new text value.
$ws.Range("G4").Formula = ""
$ws.Range("G4").Value = '$[D4 + ( E4 * ${employee.hourlyPay} )]'

# G5 = SUM(G4) naturally recalculates from #VALUE! to 0 now that G4 is no
# longer an error cell.

# Move/update the active selection to G4, matching the saved view state.
$ws.Range("G4").Select()
